$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells keep their text representation,
# since Excel would otherwise auto-convert numeric-looking strings
# like "0.999" or "1.00" into actual numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "63.884.03"
$ws.Range("E2").Value = "  -4.27%  "

# Row 3
$ws.Range("D3").Value = "3.134.40"
$ws.Range("E3").Value = "  -9.14%  "

# Row 4
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").Value = "561.01"
$ws.Range("E5").Value = "  -4.38%  "

# Row 6
$ws.Range("D6").Value = "167.94"
$ws.Range("E6").Value = "  -6.23%  "

# Row 7
$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D7").Value = "0.606"
$ws.Range("E7").Value = "  -3.83%  "

# Row 8
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.01%  "

# Row 9
$ws.Range("D9").Value = "3.129.76"
$ws.Range("E9").Value = "  -9.21%  "

# Row 10
$ws.Range("E10").Value = "  -7.87%  "

# Row 11
$ws.Range("D11").Value = "6.54"
$ws.Range("E11").Value = "  -6.09%  "

# Row 12
$ws.Range("E12").Value = "  -6.59%  "

# Row 13
$ws.Range("D13").Value = "3.664.10"
$ws.Range("E13").Value = "  -9.47%  "

# Row 14
$ws.Range("E14").Value = "  +0.46%  "

# Row 15
$ws.Range("D15").Value = "26.94"
$ws.Range("E15").Value = "  -9.76%  "

# Row 16
$ws.Range("D16").Value = "63.742.52"
$ws.Range("E16").Value = "  -4.43%  "

# Row 17
$ws.Range("E17").Value = "  -7.24%  "

# Row 18
$ws.Range("D18").Value = "3.132.03"
$ws.Range("E18").Value = "  -9.27%  "

# Row 19
$ws.Range("E19").Value = "  -4.64%  "

# Row 20
$ws.Range("D20").Value = "12.84"
$ws.Range("E20").Value = "  -7.53%  "

# Row 21
$ws.Range("D21").Value = "350.82"
$ws.Range("E21").Value = "  -5.97%  "

# Row 22
$ws.Range("D22").Value = "7.15"
$ws.Range("E22").Value = "  -6.92%  "

# Row 23
$ws.Range("D23").Value = "0.996"
$ws.Range("E23").Value = "  -0.40%  "

# Row 24
$ws.Range("D24").Value = "67.73"
$ws.Range("E24").Value = "  -7.96%  "

# Row 25
$ws.Range("D25").Value = "0.497"
$ws.Range("E25").Value = "  -7.54%  "

# Row 26
$ws.Range("D26").Value = "0.0000115"
$ws.Range("E26").Value = "  -11.14%  "

# Row 27
$ws.Range("D27").Value = "9.50"
$ws.Range("E27").Value = "  -5.03%  "

# Row 28
$ws.Range("D28").Value = "0.174"
$ws.Range("E28").Value = "  -2.39%  "

# Row 29
$ws.Range("E29").Value = "  -0.06%  "

# Row 30
$ws.Range("E30").Value = "  -0.15%  "

# Row 31
$ws.Range("D31").Value = "1.88"
$ws.Range("E31").Value = "  -6.23%  "

# Row 32
$ws.Range("D32").Value = "5.42"
$ws.Range("E32").Value = "  -8.18%  "

# Row 33
$ws.Range("D33").Value = "21.80"
$ws.Range("E33").Value = "  -7.99%  "

# Row 34
$ws.Range("E34").Value = "  -7.23%  "

# Row 35
$ws.Range("D35").Value = "6.56"
$ws.Range("E35").Value = "  -7.68%  "

# Row 36
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "1.42"
$ws.Range("E36").Value = "  -10.25%  "

# Row 37
$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").Value = "153.42"
$ws.Range("E37").Value = "  -6.05%  "

# Row 38
$ws.Range("D38").Value = "0.814"
$ws.Range("E38").Value = "  -7.91%  "

# Row 39
$ws.Range("D39").Value = "26.23"
$ws.Range("E39").Value = "  -6.19%  "

# Row 40
$ws.Range("E40").Value = "  -7.87%  "

# Row 41
$ws.Range("D41").Value = "2.633.69"
$ws.Range("E41").Value = "  -4.22%  "

# Row 42
$ws.Range("E42").Value = "  -8.32%  "

# Row 43
$ws.Range("D43").Value = "4.14"
$ws.Range("E43").Value = "  -8.33%  "

# Row 44
$ws.Range("D44").Value = "39.10"
$ws.Range("E44").Value = "  -2.51%  "

# Row 45
$ws.Range("D45").Value = "5.93"
$ws.Range("E45").Value = "  -7.81%  "

# Row 46
$ws.Range("E46").Value = "  -7.91%  "

# Row 47
$ws.Range("D47").Value = "23.57"
$ws.Range("E47").Value = "  -7.95%  "

# Row 48
$ws.Range("D48").Value = "313.97"
$ws.Range("E48").Value = "  -7.06%  "

# Row 49
$ws.Range("D49").Value = "0.0269"
$ws.Range("E49").Value = "  -6.85%  "

# Row 50
$ws.Range("E50").Value = "  -4.71%  "

# Row 51
$ws.Range("D51").Value = "0.998"
$ws.Range("E51").Value = "  -0.09%  "
